$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.970.99"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "2.343.21"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'541.59"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'133.98"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("D8").Value = "'0.563"
$ws.Range("E8").Value = "  +4.84%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "'5.52"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "'23.79"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").Value = "2.760.38"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "57.916.93"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "'0.0000135"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "2.343.34"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "'10.67"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "'328.35"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'62.85"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("D25").Value = "'0.996"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'8.30"
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("E27").Value = "  -6.41%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "'170.34"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "0.0₃0734"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").Value = "'6.13"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("D36").Value = "'4.15"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  -2.34%  "
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").Value = "'141.36"
$ws.Range("E40").Value = "  -6.03%  "
$ws.Range("D41").Value = "'0.377"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").Value = "'288.33"
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").Value = "'3.63"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "'0.0945"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").Value = "'18.98"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").Value = "'0.378"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").Value = "'11.10"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("E51").Value = "  +0.83%  "
